# Fix format upload file: replace the product-upload template headers
# (qty/nama_barang/harga/harga_jual/jenis/kategori) with the
# people/registration-upload template headers
# (nama/tgl_lahir/jenis_kelamin/alamat/sektor), drop the now-unused 6th
# column, add an example date cell formatted as a date below tgl_lahir,
# widen the jenis_kelamin column, and leave the selection on B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:E1) -------------------------------------------------
# Set the cells in left-to-right order so the shared-strings table is
# rebuilt in the same order the strings are first used.
$ws.Range("A1").Value = "nama"
$ws.Range("B1").Value = "tgl_lahir"
$ws.Range("C1").Value = "jenis_kelamin"
$ws.Range("D1").Value = "alamat"
$ws.Range("E1").Value = "sektor"

# The old template had a 6th column (harga_jual) that is no longer used.
$ws.Range("F1").ClearContents()

# --- Example date cell under tgl_lahir ----------------------------------
$ws.Range("B2").NumberFormat = "mm-dd-yy"

# --- Column width for jenis_kelamin -------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.498697916666668

# --- Selection ------------------------------------------------------------
$ws.Range("B3").Select()
